# Auto-generated edit script: relocates every data row (2-115) of the
# "Avverkningsanmälningar" sheet to its new position and bumps the
# "Förändrad" (column C) timestamp for every data row to 46066.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 2  <=  old row 2  (A 46651-2024)
$ws.Range("A2:Z2").ClearContents()
$ws.Range("A2").Value = 'A 46651-2024'
$ws.Range("B2").Value = 45583
$ws.Range("D2").Value = 'STOCKHOLMS LÄN'
$ws.Range("E2").Value = 'ÖSTERÅKER'
$ws.Range("G2").Value = 8.1
$ws.Range("H2").Value = 5
$ws.Range("I2").Value = 15
$ws.Range("J2").Value = 8
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 11
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 28
$ws.Range("R2").Value = "Blackticka`r`nKnärot`r`nVågticka`r`nGarnlav`r`nGrön aspvedbock`r`nKådvaxskinn`r`nLunglav`r`nSpillkråka`r`nTallticka`r`nUllticka`r`nVedtrappmossa`r`nBjörksplintborre`r`nBronshjon`r`nFlagellkvastmossa`r`nGranbarkgnagare`r`nGrön sköldmossa`r`nLönnlav`r`nMindre märgborre`r`nMörk husmossa`r`nRostfläck`r`nSkarp dropptaggsvamp`r`nSkinnlav`r`nStor aspticka`r`nStrutbräken`r`nVedticka`r`nVågbandad barkbock`r`nBlåsippa`r`nRevlummer"
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/artfynd/A 46651-2024 artfynd.xlsx", "A 46651-2024")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/kartor/A 46651-2024 karta.png", "A 46651-2024")'
$ws.Range("U2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/knärot/A 46651-2024 karta knärot.png", "A 46651-2024")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomål/A 46651-2024 FSC-klagomål.docx", "A 46651-2024")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomålsmail/A 46651-2024 FSC-klagomål mail.docx", "A 46651-2024")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsyn/A 46651-2024 tillsynsbegäran.docx", "A 46651-2024")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsynsmail/A 46651-2024 tillsynsbegäran mail.docx", "A 46651-2024")'
$ws.Range("Z2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/fåglar/A 46651-2024 prioriterade fågelarter.docx", "A 46651-2024")'
$ws.Range("C2").Value = 46066

# New row 3  <=  old row 3  (A 31459-2025)
$ws.Range("A3:Z3").ClearContents()
$ws.Range("A3").Value = 'A 31459-2025'
$ws.Range("B3").Value = 45833
$ws.Range("D3").Value = 'STOCKHOLMS LÄN'
$ws.Range("E3").Value = 'ÖSTERÅKER'
$ws.Range("G3").Value = 5.7
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = 3
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 3
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 8
$ws.Range("R3").Value = "Garnlav`r`nTalltita`r`nUllticka`r`nMörk husmossa`r`nSårläka`r`nThomsons trägnagare`r`nNattviol`r`nBlåsippa"
$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/artfynd/A 31459-2025 artfynd.xlsx", "A 31459-2025")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/kartor/A 31459-2025 karta.png", "A 31459-2025")'
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomål/A 31459-2025 FSC-klagomål.docx", "A 31459-2025")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomålsmail/A 31459-2025 FSC-klagomål mail.docx", "A 31459-2025")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsyn/A 31459-2025 tillsynsbegäran.docx", "A 31459-2025")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsynsmail/A 31459-2025 tillsynsbegäran mail.docx", "A 31459-2025")'
$ws.Range("Z3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/fåglar/A 31459-2025 prioriterade fågelarter.docx", "A 31459-2025")'
$ws.Range("C3").Value = 46066

# New row 4  <=  old row 4  (A 31015-2024)
$ws.Range("A4:Z4").ClearContents()
$ws.Range("A4").Value = 'A 31015-2024'
$ws.Range("B4").Value = 45499
$ws.Range("D4").Value = 'STOCKHOLMS LÄN'
$ws.Range("E4").Value = 'ÖSTERÅKER'
$ws.Range("F4").Value = 'Övriga Aktiebolag'
$ws.Range("G4").Value = 5.3
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 4
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 3
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 7
$ws.Range("R4").Value = "Spillkråka`r`nTallticka`r`nUllticka`r`nBlodticka`r`nBronshjon`r`nJättesvampmal`r`nVågbandad barkbock"
$ws.Range("S4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/artfynd/A 31015-2024 artfynd.xlsx", "A 31015-2024")'
$ws.Range("T4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/kartor/A 31015-2024 karta.png", "A 31015-2024")'
$ws.Range("V4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomål/A 31015-2024 FSC-klagomål.docx", "A 31015-2024")'
$ws.Range("W4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomålsmail/A 31015-2024 FSC-klagomål mail.docx", "A 31015-2024")'
$ws.Range("X4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsyn/A 31015-2024 tillsynsbegäran.docx", "A 31015-2024")'
$ws.Range("Y4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsynsmail/A 31015-2024 tillsynsbegäran mail.docx", "A 31015-2024")'
$ws.Range("Z4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/fåglar/A 31015-2024 prioriterade fågelarter.docx", "A 31015-2024")'
$ws.Range("C4").Value = 46066

# New row 5  <=  old row 5  (A 62573-2025)
$ws.Range("A5:Z5").ClearContents()
$ws.Range("A5").Value = 'A 62573-2025'
$ws.Range("B5").Value = 46007.73069444444
$ws.Range("D5").Value = 'STOCKHOLMS LÄN'
$ws.Range("E5").Value = 'ÖSTERÅKER'
$ws.Range("G5").Value = 3.9
$ws.Range("H5").Value = 5
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 4
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 4
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 7
$ws.Range("R5").Value = "Nordfladdermus`r`nSlåtterfibbla`r`nSvartvit flugsnappare`r`nÄrtsångare`r`nVårärt`r`nStörre brunfladdermus`r`nBlåsippa"
$ws.Range("S5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/artfynd/A 62573-2025 artfynd.xlsx", "A 62573-2025")'
$ws.Range("T5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/kartor/A 62573-2025 karta.png", "A 62573-2025")'
$ws.Range("V5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomål/A 62573-2025 FSC-klagomål.docx", "A 62573-2025")'
$ws.Range("W5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomålsmail/A 62573-2025 FSC-klagomål mail.docx", "A 62573-2025")'
$ws.Range("X5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsyn/A 62573-2025 tillsynsbegäran.docx", "A 62573-2025")'
$ws.Range("Y5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsynsmail/A 62573-2025 tillsynsbegäran mail.docx", "A 62573-2025")'
$ws.Range("Z5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/fåglar/A 62573-2025 prioriterade fågelarter.docx", "A 62573-2025")'
$ws.Range("C5").Value = 46066

# New row 6  <=  old row 6  (A 62642-2025)
$ws.Range("A6:Z6").ClearContents()
$ws.Range("A6").Value = 'A 62642-2025'
$ws.Range("B6").Value = 46008.37217592593
$ws.Range("D6").Value = 'STOCKHOLMS LÄN'
$ws.Range("E6").Value = 'ÖSTERÅKER'
$ws.Range("G6").Value = 2.7
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 6
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 1
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 7
$ws.Range("R6").Value = "Grön aspvedbock`r`nBjörksplintborre`r`nBronshjon`r`nDvärgtufs`r`nGranbarkgnagare`r`nMindre märgborre`r`nVågbandad barkbock"
$ws.Range("S6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/artfynd/A 62642-2025 artfynd.xlsx", "A 62642-2025")'
$ws.Range("T6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/kartor/A 62642-2025 karta.png", "A 62642-2025")'
$ws.Range("V6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomål/A 62642-2025 FSC-klagomål.docx", "A 62642-2025")'
$ws.Range("W6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomålsmail/A 62642-2025 FSC-klagomål mail.docx", "A 62642-2025")'
$ws.Range("X6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsyn/A 62642-2025 tillsynsbegäran.docx", "A 62642-2025")'
$ws.Range("Y6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsynsmail/A 62642-2025 tillsynsbegäran mail.docx", "A 62642-2025")'
$ws.Range("C6").Value = 46066

# New row 7  <=  old row 7  (A 153-2025)
$ws.Range("A7:Z7").ClearContents()
$ws.Range("A7").Value = 'A 153-2025'
$ws.Range("B7").Value = 45659.61650462963
$ws.Range("D7").Value = 'STOCKHOLMS LÄN'
$ws.Range("E7").Value = 'ÖSTERÅKER'
$ws.Range("F7").Value = 'Kommuner'
$ws.Range("G7").Value = 5.3
$ws.Range("H7").Value = 5
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 1
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 4
$ws.Range("P7").Value = 3
$ws.Range("Q7").Value = 5
$ws.Range("R7").Value = "Grönfink`r`nStare`r`nTallbit`r`nSvartvit flugsnappare`r`nGråkråka"
$ws.Range("S7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/artfynd/A 153-2025 artfynd.xlsx", "A 153-2025")'
$ws.Range("T7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/kartor/A 153-2025 karta.png", "A 153-2025")'
$ws.Range("V7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomål/A 153-2025 FSC-klagomål.docx", "A 153-2025")'
$ws.Range("W7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomålsmail/A 153-2025 FSC-klagomål mail.docx", "A 153-2025")'
$ws.Range("X7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsyn/A 153-2025 tillsynsbegäran.docx", "A 153-2025")'
$ws.Range("Y7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsynsmail/A 153-2025 tillsynsbegäran mail.docx", "A 153-2025")'
$ws.Range("Z7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/fåglar/A 153-2025 prioriterade fågelarter.docx", "A 153-2025")'
$ws.Range("C7").Value = 46066

# New row 8  <=  old row 9  (A 34717-2025)
$ws.Range("A8:Z8").ClearContents()
$ws.Range("A8").Value = 'A 34717-2025'
$ws.Range("B8").Value = 45848
$ws.Range("D8").Value = 'STOCKHOLMS LÄN'
$ws.Range("E8").Value = 'ÖSTERÅKER'
$ws.Range("G8").Value = 3.8
$ws.Range("H8").Value = 3
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 1
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 0
$ws.Range("O8").Value = 1
$ws.Range("P8").Value = 0
$ws.Range("Q8").Value = 3
$ws.Range("R8").Value = "Nordfladdermus`r`nDvärgpipistrell`r`nStörre brunfladdermus"
$ws.Range("S8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/artfynd/A 34717-2025 artfynd.xlsx", "A 34717-2025")'
$ws.Range("T8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/kartor/A 34717-2025 karta.png", "A 34717-2025")'
$ws.Range("V8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomål/A 34717-2025 FSC-klagomål.docx", "A 34717-2025")'
$ws.Range("W8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomålsmail/A 34717-2025 FSC-klagomål mail.docx", "A 34717-2025")'
$ws.Range("X8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsyn/A 34717-2025 tillsynsbegäran.docx", "A 34717-2025")'
$ws.Range("Y8").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsynsmail/A 34717-2025 tillsynsbegäran mail.docx", "A 34717-2025")'
$ws.Range("C8").Value = 46066

# New row 9  <=  old row 8  (A 12620-2022)
$ws.Range("A9:Z9").ClearContents()
$ws.Range("A9").Value = 'A 12620-2022'
$ws.Range("B9").Value = 44638
$ws.Range("D9").Value = 'STOCKHOLMS LÄN'
$ws.Range("E9").Value = 'ÖSTERÅKER'
$ws.Range("G9").Value = 9.4
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 0
$ws.Range("O9").Value = 1
$ws.Range("P9").Value = 1
$ws.Range("Q9").Value = 3
$ws.Range("R9").Value = "Taggfingersvamp`r`nFällmossa`r`nBlåsippa"
$ws.Range("S9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/artfynd/A 12620-2022 artfynd.xlsx", "A 12620-2022")'
$ws.Range("T9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/kartor/A 12620-2022 karta.png", "A 12620-2022")'
$ws.Range("V9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomål/A 12620-2022 FSC-klagomål.docx", "A 12620-2022")'
$ws.Range("W9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomålsmail/A 12620-2022 FSC-klagomål mail.docx", "A 12620-2022")'
$ws.Range("X9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsyn/A 12620-2022 tillsynsbegäran.docx", "A 12620-2022")'
$ws.Range("Y9").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsynsmail/A 12620-2022 tillsynsbegäran mail.docx", "A 12620-2022")'
$ws.Range("C9").Value = 46066

# New row 10  <=  old row 10  (A 17244-2022)
$ws.Range("A10:Z10").ClearContents()
$ws.Range("A10").Value = 'A 17244-2022'
$ws.Range("B10").Value = 44678
$ws.Range("D10").Value = 'STOCKHOLMS LÄN'
$ws.Range("E10").Value = 'ÖSTERÅKER'
$ws.Range("G10").Value = 1.4
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 2
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 2
$ws.Range("R10").Value = "Bronshjon`r`nGrön sköldmossa"
$ws.Range("S10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/artfynd/A 17244-2022 artfynd.xlsx", "A 17244-2022")'
$ws.Range("T10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/kartor/A 17244-2022 karta.png", "A 17244-2022")'
$ws.Range("V10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomål/A 17244-2022 FSC-klagomål.docx", "A 17244-2022")'
$ws.Range("W10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomålsmail/A 17244-2022 FSC-klagomål mail.docx", "A 17244-2022")'
$ws.Range("X10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsyn/A 17244-2022 tillsynsbegäran.docx", "A 17244-2022")'
$ws.Range("Y10").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsynsmail/A 17244-2022 tillsynsbegäran mail.docx", "A 17244-2022")'
$ws.Range("C10").Value = 46066

# New row 11  <=  old row 11  (A 61275-2021)
$ws.Range("A11:Z11").ClearContents()
$ws.Range("A11").Value = 'A 61275-2021'
$ws.Range("B11").Value = 44498
$ws.Range("D11").Value = 'STOCKHOLMS LÄN'
$ws.Range("E11").Value = 'ÖSTERÅKER'
$ws.Range("G11").Value = 2.1
$ws.Range("H11").Value = 2
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 0
$ws.Range("Q11").Value = 2
$ws.Range("R11").Value = "Fläcknycklar`r`nNattviol"
$ws.Range("S11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/artfynd/A 61275-2021 artfynd.xlsx", "A 61275-2021")'
$ws.Range("T11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/kartor/A 61275-2021 karta.png", "A 61275-2021")'
$ws.Range("V11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomål/A 61275-2021 FSC-klagomål.docx", "A 61275-2021")'
$ws.Range("W11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomålsmail/A 61275-2021 FSC-klagomål mail.docx", "A 61275-2021")'
$ws.Range("X11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsyn/A 61275-2021 tillsynsbegäran.docx", "A 61275-2021")'
$ws.Range("Y11").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsynsmail/A 61275-2021 tillsynsbegäran mail.docx", "A 61275-2021")'
$ws.Range("C11").Value = 46066

# New row 12  <=  old row 13  (A 36366-2025)
$ws.Range("A12:Z12").ClearContents()
$ws.Range("A12").Value = 'A 36366-2025'
$ws.Range("B12").Value = 45868
$ws.Range("D12").Value = 'STOCKHOLMS LÄN'
$ws.Range("E12").Value = 'ÖSTERÅKER'
$ws.Range("G12").Value = 7.2
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 1
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 0
$ws.Range("N12").Value = 0
$ws.Range("O12").Value = 1
$ws.Range("P12").Value = 0
$ws.Range("Q12").Value = 2
$ws.Range("R12").Value = "Garnlav`r`nBjörksplintborre"
$ws.Range("S12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/artfynd/A 36366-2025 artfynd.xlsx", "A 36366-2025")'
$ws.Range("T12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/kartor/A 36366-2025 karta.png", "A 36366-2025")'
$ws.Range("V12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomål/A 36366-2025 FSC-klagomål.docx", "A 36366-2025")'
$ws.Range("W12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomålsmail/A 36366-2025 FSC-klagomål mail.docx", "A 36366-2025")'
$ws.Range("X12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsyn/A 36366-2025 tillsynsbegäran.docx", "A 36366-2025")'
$ws.Range("Y12").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsynsmail/A 36366-2025 tillsynsbegäran mail.docx", "A 36366-2025")'
$ws.Range("C12").Value = 46066

# New row 13  <=  old row 17  (A 34721-2025)
$ws.Range("A13:Z13").ClearContents()
$ws.Range("A13").Value = 'A 34721-2025'
$ws.Range("B13").Value = 45848
$ws.Range("D13").Value = 'STOCKHOLMS LÄN'
$ws.Range("E13").Value = 'ÖSTERÅKER'
$ws.Range("G13").Value = 2.6
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 2
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = 0
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 0
$ws.Range("Q13").Value = 2
$ws.Range("R13").Value = "Skuggblåslav`r`nVästlig hakmossa"
$ws.Range("S13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/artfynd/A 34721-2025 artfynd.xlsx", "A 34721-2025")'
$ws.Range("T13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/kartor/A 34721-2025 karta.png", "A 34721-2025")'
$ws.Range("V13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomål/A 34721-2025 FSC-klagomål.docx", "A 34721-2025")'
$ws.Range("W13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomålsmail/A 34721-2025 FSC-klagomål mail.docx", "A 34721-2025")'
$ws.Range("X13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsyn/A 34721-2025 tillsynsbegäran.docx", "A 34721-2025")'
$ws.Range("Y13").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsynsmail/A 34721-2025 tillsynsbegäran mail.docx", "A 34721-2025")'
$ws.Range("C13").Value = 46066

# New row 14  <=  old row 14  (A 43308-2024)
$ws.Range("A14:Z14").ClearContents()
$ws.Range("A14").Value = 'A 43308-2024'
$ws.Range("B14").Value = 45568
$ws.Range("D14").Value = 'STOCKHOLMS LÄN'
$ws.Range("E14").Value = 'ÖSTERÅKER'
$ws.Range("F14").Value = 'Kyrkan'
$ws.Range("G14").Value = 3.2
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 1
$ws.Range("J14").Value = 1
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 0
$ws.Range("O14").Value = 1
$ws.Range("P14").Value = 0
$ws.Range("Q14").Value = 2
$ws.Range("R14").Value = "Tallticka`r`nDropptaggsvamp"
$ws.Range("S14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/artfynd/A 43308-2024 artfynd.xlsx", "A 43308-2024")'
$ws.Range("T14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/kartor/A 43308-2024 karta.png", "A 43308-2024")'
$ws.Range("V14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomål/A 43308-2024 FSC-klagomål.docx", "A 43308-2024")'
$ws.Range("W14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomålsmail/A 43308-2024 FSC-klagomål mail.docx", "A 43308-2024")'
$ws.Range("X14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsyn/A 43308-2024 tillsynsbegäran.docx", "A 43308-2024")'
$ws.Range("Y14").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsynsmail/A 43308-2024 tillsynsbegäran mail.docx", "A 43308-2024")'
$ws.Range("C14").Value = 46066

# New row 15  <=  old row 18  (A 62578-2025)
$ws.Range("A15:Z15").ClearContents()
$ws.Range("A15").Value = 'A 62578-2025'
$ws.Range("B15").Value = 46007
$ws.Range("D15").Value = 'STOCKHOLMS LÄN'
$ws.Range("E15").Value = 'ÖSTERÅKER'
$ws.Range("G15").Value = 3.4
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 2
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = 0
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = 0
$ws.Range("Q15").Value = 2
$ws.Range("R15").Value = "Fällmossa`r`nMindre märgborre"
$ws.Range("S15").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/artfynd/A 62578-2025 artfynd.xlsx", "A 62578-2025")'
$ws.Range("T15").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/kartor/A 62578-2025 karta.png", "A 62578-2025")'
$ws.Range("V15").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomål/A 62578-2025 FSC-klagomål.docx", "A 62578-2025")'
$ws.Range("W15").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomålsmail/A 62578-2025 FSC-klagomål mail.docx", "A 62578-2025")'
$ws.Range("X15").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsyn/A 62578-2025 tillsynsbegäran.docx", "A 62578-2025")'
$ws.Range("Y15").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsynsmail/A 62578-2025 tillsynsbegäran mail.docx", "A 62578-2025")'
$ws.Range("C15").Value = 46066

# New row 16  <=  old row 12  (A 31523-2022)
$ws.Range("A16:Z16").ClearContents()
$ws.Range("A16").Value = 'A 31523-2022'
$ws.Range("B16").Value = 44775.22490740741
$ws.Range("D16").Value = 'STOCKHOLMS LÄN'
$ws.Range("E16").Value = 'ÖSTERÅKER'
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 2
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = 0
$ws.Range("O16").Value = 2
$ws.Range("P16").Value = 0
$ws.Range("Q16").Value = 2
$ws.Range("R16").Value = "Garnlav`r`nTallticka"
$ws.Range("S16").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/artfynd/A 31523-2022 artfynd.xlsx", "A 31523-2022")'
$ws.Range("T16").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/kartor/A 31523-2022 karta.png", "A 31523-2022")'
$ws.Range("V16").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomål/A 31523-2022 FSC-klagomål.docx", "A 31523-2022")'
$ws.Range("W16").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomålsmail/A 31523-2022 FSC-klagomål mail.docx", "A 31523-2022")'
$ws.Range("X16").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsyn/A 31523-2022 tillsynsbegäran.docx", "A 31523-2022")'
$ws.Range("Y16").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsynsmail/A 31523-2022 tillsynsbegäran mail.docx", "A 31523-2022")'
$ws.Range("C16").Value = 46066

# New row 17  <=  old row 16  (A 61363-2022)
$ws.Range("A17:Z17").ClearContents()
$ws.Range("A17").Value = 'A 61363-2022'
$ws.Range("B17").Value = 44915
$ws.Range("D17").Value = 'STOCKHOLMS LÄN'
$ws.Range("E17").Value = 'ÖSTERÅKER'
$ws.Range("G17").Value = 5.3
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 1
$ws.Range("J17").Value = 1
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = 0
$ws.Range("O17").Value = 1
$ws.Range("P17").Value = 0
$ws.Range("Q17").Value = 2
$ws.Range("R17").Value = "Motaggsvamp`r`nSkuggblåslav"
$ws.Range("S17").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/artfynd/A 61363-2022 artfynd.xlsx", "A 61363-2022")'
$ws.Range("T17").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/kartor/A 61363-2022 karta.png", "A 61363-2022")'
$ws.Range("V17").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomål/A 61363-2022 FSC-klagomål.docx", "A 61363-2022")'
$ws.Range("W17").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomålsmail/A 61363-2022 FSC-klagomål mail.docx", "A 61363-2022")'
$ws.Range("X17").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsyn/A 61363-2022 tillsynsbegäran.docx", "A 61363-2022")'
$ws.Range("Y17").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsynsmail/A 61363-2022 tillsynsbegäran mail.docx", "A 61363-2022")'
$ws.Range("C17").Value = 46066

# New row 18  <=  old row 15  (A 54703-2024)
$ws.Range("A18:Z18").ClearContents()
$ws.Range("A18").Value = 'A 54703-2024'
$ws.Range("B18").Value = 45618
$ws.Range("D18").Value = 'STOCKHOLMS LÄN'
$ws.Range("E18").Value = 'ÖSTERÅKER'
$ws.Range("G18").Value = 7.7
$ws.Range("H18").Value = 1
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 2
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 0
$ws.Range("N18").Value = 0
$ws.Range("O18").Value = 2
$ws.Range("P18").Value = 0
$ws.Range("Q18").Value = 2
$ws.Range("R18").Value = "Spillkråka`r`nTallticka"
$ws.Range("S18").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/artfynd/A 54703-2024 artfynd.xlsx", "A 54703-2024")'
$ws.Range("T18").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/kartor/A 54703-2024 karta.png", "A 54703-2024")'
$ws.Range("V18").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomål/A 54703-2024 FSC-klagomål.docx", "A 54703-2024")'
$ws.Range("W18").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomålsmail/A 54703-2024 FSC-klagomål mail.docx", "A 54703-2024")'
$ws.Range("X18").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsyn/A 54703-2024 tillsynsbegäran.docx", "A 54703-2024")'
$ws.Range("Y18").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsynsmail/A 54703-2024 tillsynsbegäran mail.docx", "A 54703-2024")'
$ws.Range("Z18").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/fåglar/A 54703-2024 prioriterade fågelarter.docx", "A 54703-2024")'
$ws.Range("C18").Value = 46066

# New row 19  <=  old row 19  (A 22172-2022)
$ws.Range("A19:Z19").ClearContents()
$ws.Range("A19").Value = 'A 22172-2022'
$ws.Range("B19").Value = 44712
$ws.Range("D19").Value = 'STOCKHOLMS LÄN'
$ws.Range("E19").Value = 'ÖSTERÅKER'
$ws.Range("G19").Value = 3.7
$ws.Range("H19").Value = 1
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 1
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = 0
$ws.Range("O19").Value = 1
$ws.Range("P19").Value = 0
$ws.Range("Q19").Value = 1
$ws.Range("R19").Value = 'Kråka'
$ws.Range("S19").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/artfynd/A 22172-2022 artfynd.xlsx", "A 22172-2022")'
$ws.Range("T19").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/kartor/A 22172-2022 karta.png", "A 22172-2022")'
$ws.Range("V19").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomål/A 22172-2022 FSC-klagomål.docx", "A 22172-2022")'
$ws.Range("W19").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomålsmail/A 22172-2022 FSC-klagomål mail.docx", "A 22172-2022")'
$ws.Range("X19").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsyn/A 22172-2022 tillsynsbegäran.docx", "A 22172-2022")'
$ws.Range("Y19").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsynsmail/A 22172-2022 tillsynsbegäran mail.docx", "A 22172-2022")'
$ws.Range("Z19").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/fåglar/A 22172-2022 prioriterade fågelarter.docx", "A 22172-2022")'
$ws.Range("C19").Value = 46066

# New row 20  <=  old row 21  (A 27921-2025)
$ws.Range("A20:Z20").ClearContents()
$ws.Range("A20").Value = 'A 27921-2025'
$ws.Range("B20").Value = 45817.4846875
$ws.Range("D20").Value = 'STOCKHOLMS LÄN'
$ws.Range("E20").Value = 'ÖSTERÅKER'
$ws.Range("G20").Value = 11.8
$ws.Range("H20").Value = 1
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = 0
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = 0
$ws.Range("Q20").Value = 1
$ws.Range("R20").Value = 'Mindre vattensalamander'
$ws.Range("S20").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/artfynd/A 27921-2025 artfynd.xlsx", "A 27921-2025")'
$ws.Range("T20").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/kartor/A 27921-2025 karta.png", "A 27921-2025")'
$ws.Range("V20").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomål/A 27921-2025 FSC-klagomål.docx", "A 27921-2025")'
$ws.Range("W20").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomålsmail/A 27921-2025 FSC-klagomål mail.docx", "A 27921-2025")'
$ws.Range("X20").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsyn/A 27921-2025 tillsynsbegäran.docx", "A 27921-2025")'
$ws.Range("Y20").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsynsmail/A 27921-2025 tillsynsbegäran mail.docx", "A 27921-2025")'
$ws.Range("C20").Value = 46066

# New row 21  <=  old row 23  (A 31716-2025)
$ws.Range("A21:Z21").ClearContents()
$ws.Range("A21").Value = 'A 31716-2025'
$ws.Range("B21").Value = 45834
$ws.Range("D21").Value = 'STOCKHOLMS LÄN'
$ws.Range("E21").Value = 'ÖSTERÅKER'
$ws.Range("G21").Value = 2.8
$ws.Range("H21").Value = 1
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 0
$ws.Range("N21").Value = 0
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = 0
$ws.Range("Q21").Value = 1
$ws.Range("R21").Value = 'Revlummer'
$ws.Range("S21").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/artfynd/A 31716-2025 artfynd.xlsx", "A 31716-2025")'
$ws.Range("T21").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/kartor/A 31716-2025 karta.png", "A 31716-2025")'
$ws.Range("V21").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomål/A 31716-2025 FSC-klagomål.docx", "A 31716-2025")'
$ws.Range("W21").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomålsmail/A 31716-2025 FSC-klagomål mail.docx", "A 31716-2025")'
$ws.Range("X21").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsyn/A 31716-2025 tillsynsbegäran.docx", "A 31716-2025")'
$ws.Range("Y21").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsynsmail/A 31716-2025 tillsynsbegäran mail.docx", "A 31716-2025")'
$ws.Range("C21").Value = 46066

# New row 22  <=  old row 24  (A 38968-2025)
$ws.Range("A22:Z22").ClearContents()
$ws.Range("A22").Value = 'A 38968-2025'
$ws.Range("B22").Value = 45887.70863425926
$ws.Range("D22").Value = 'STOCKHOLMS LÄN'
$ws.Range("E22").Value = 'ÖSTERÅKER'
$ws.Range("G22").Value = 7.2
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 1
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 0
$ws.Range("N22").Value = 0
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = 0
$ws.Range("Q22").Value = 1
$ws.Range("R22").Value = 'Fjällig taggsvamp s.str.'
$ws.Range("S22").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/artfynd/A 38968-2025 artfynd.xlsx", "A 38968-2025")'
$ws.Range("T22").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/kartor/A 38968-2025 karta.png", "A 38968-2025")'
$ws.Range("V22").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomål/A 38968-2025 FSC-klagomål.docx", "A 38968-2025")'
$ws.Range("W22").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomålsmail/A 38968-2025 FSC-klagomål mail.docx", "A 38968-2025")'
$ws.Range("X22").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsyn/A 38968-2025 tillsynsbegäran.docx", "A 38968-2025")'
$ws.Range("Y22").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsynsmail/A 38968-2025 tillsynsbegäran mail.docx", "A 38968-2025")'
$ws.Range("C22").Value = 46066

# New row 23  <=  old row 25  (A 3914-2023)
$ws.Range("A23:Z23").ClearContents()
$ws.Range("A23").Value = 'A 3914-2023'
$ws.Range("B23").Value = 44952.30934027778
$ws.Range("D23").Value = 'STOCKHOLMS LÄN'
$ws.Range("E23").Value = 'ÖSTERÅKER'
$ws.Range("F23").Value = 'Kommuner'
$ws.Range("G23").Value = 15.3
$ws.Range("H23").Value = 1
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 0
$ws.Range("N23").Value = 0
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = 0
$ws.Range("Q23").Value = 1
$ws.Range("R23").Value = 'Gråkråka'
$ws.Range("S23").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/artfynd/A 3914-2023 artfynd.xlsx", "A 3914-2023")'
$ws.Range("T23").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/kartor/A 3914-2023 karta.png", "A 3914-2023")'
$ws.Range("V23").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomål/A 3914-2023 FSC-klagomål.docx", "A 3914-2023")'
$ws.Range("W23").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomålsmail/A 3914-2023 FSC-klagomål mail.docx", "A 3914-2023")'
$ws.Range("X23").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsyn/A 3914-2023 tillsynsbegäran.docx", "A 3914-2023")'
$ws.Range("Y23").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsynsmail/A 3914-2023 tillsynsbegäran mail.docx", "A 3914-2023")'
$ws.Range("C23").Value = 46066

# New row 24  <=  old row 20  (A 13942-2021)
$ws.Range("A24:Z24").ClearContents()
$ws.Range("A24").Value = 'A 13942-2021'
$ws.Range("B24").Value = 44277
$ws.Range("D24").Value = 'STOCKHOLMS LÄN'
$ws.Range("E24").Value = 'ÖSTERÅKER'
$ws.Range("G24").Value = 5.2
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 1
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 0
$ws.Range("N24").Value = 0
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = 0
$ws.Range("Q24").Value = 1
$ws.Range("R24").Value = 'Skogsbräsma'
$ws.Range("S24").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/artfynd/A 13942-2021 artfynd.xlsx", "A 13942-2021")'
$ws.Range("T24").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/kartor/A 13942-2021 karta.png", "A 13942-2021")'
$ws.Range("V24").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomål/A 13942-2021 FSC-klagomål.docx", "A 13942-2021")'
$ws.Range("W24").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomålsmail/A 13942-2021 FSC-klagomål mail.docx", "A 13942-2021")'
$ws.Range("X24").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsyn/A 13942-2021 tillsynsbegäran.docx", "A 13942-2021")'
$ws.Range("Y24").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsynsmail/A 13942-2021 tillsynsbegäran mail.docx", "A 13942-2021")'
$ws.Range("C24").Value = 46066

# New row 25  <=  old row 22  (A 3919-2023)
$ws.Range("A25:Z25").ClearContents()
$ws.Range("A25").Value = 'A 3919-2023'
$ws.Range("B25").Value = 44952
$ws.Range("D25").Value = 'STOCKHOLMS LÄN'
$ws.Range("E25").Value = 'ÖSTERÅKER'
$ws.Range("F25").Value = 'Kommuner'
$ws.Range("G25").Value = 4.9
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 1
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0
$ws.Range("N25").Value = 0
$ws.Range("O25").Value = 1
$ws.Range("P25").Value = 0
$ws.Range("Q25").Value = 1
$ws.Range("R25").Value = 'Tallticka'
$ws.Range("S25").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/artfynd/A 3919-2023 artfynd.xlsx", "A 3919-2023")'
$ws.Range("T25").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/kartor/A 3919-2023 karta.png", "A 3919-2023")'
$ws.Range("V25").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomål/A 3919-2023 FSC-klagomål.docx", "A 3919-2023")'
$ws.Range("W25").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/klagomålsmail/A 3919-2023 FSC-klagomål mail.docx", "A 3919-2023")'
$ws.Range("X25").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsyn/A 3919-2023 tillsynsbegäran.docx", "A 3919-2023")'
$ws.Range("Y25").Formula = '=HYPERLINK("https://klasma.github.io/Logging_0117/tillsynsmail/A 3919-2023 tillsynsbegäran mail.docx", "A 3919-2023")'
$ws.Range("C25").Value = 46066

# New row 26  <=  old row 26  (A 848-2022)
$ws.Range("A26:Z26").ClearContents()
$ws.Range("A26").Value = 'A 848-2022'
$ws.Range("B26").Value = 44569
$ws.Range("D26").Value = 'STOCKHOLMS LÄN'
$ws.Range("E26").Value = 'ÖSTERÅKER'
$ws.Range("F26").Value = 'Övriga Aktiebolag'
$ws.Range("G26").Value = 5
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = 0
$ws.Range("N26").Value = 0
$ws.Range("O26").Value = 0
$ws.Range("P26").Value = 0
$ws.Range("Q26").Value = 0
$ws.Range("R26").Value = ''
$ws.Range("C26").Value = 46066

# New row 27  <=  old row 27  (A 2148-2022)
$ws.Range("A27:Z27").ClearContents()
$ws.Range("A27").Value = 'A 2148-2022'
$ws.Range("B27").Value = 44575
$ws.Range("D27").Value = 'STOCKHOLMS LÄN'
$ws.Range("E27").Value = 'ÖSTERÅKER'
$ws.Range("F27").Value = 'Övriga Aktiebolag'
$ws.Range("G27").Value = 0.6
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = 0
$ws.Range("N27").Value = 0
$ws.Range("O27").Value = 0
$ws.Range("P27").Value = 0
$ws.Range("Q27").Value = 0
$ws.Range("R27").Value = ''
$ws.Range("C27").Value = 46066

# New row 28  <=  old row 28  (A 1910-2022)
$ws.Range("A28:Z28").ClearContents()
$ws.Range("A28").Value = 'A 1910-2022'
$ws.Range("B28").Value = 44575
$ws.Range("D28").Value = 'STOCKHOLMS LÄN'
$ws.Range("E28").Value = 'ÖSTERÅKER'
$ws.Range("F28").Value = 'Övriga Aktiebolag'
$ws.Range("G28").Value = 2.6
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 0
$ws.Range("N28").Value = 0
$ws.Range("O28").Value = 0
$ws.Range("P28").Value = 0
$ws.Range("Q28").Value = 0
$ws.Range("R28").Value = ''
$ws.Range("C28").Value = 46066

# New row 29  <=  old row 29  (A 72508-2021)
$ws.Range("A29:Z29").ClearContents()
$ws.Range("A29").Value = 'A 72508-2021'
$ws.Range("B29").Value = 44544
$ws.Range("D29").Value = 'STOCKHOLMS LÄN'
$ws.Range("E29").Value = 'ÖSTERÅKER'
$ws.Range("F29").Value = 'Övriga Aktiebolag'
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = 0
$ws.Range("O29").Value = 0
$ws.Range("P29").Value = 0
$ws.Range("Q29").Value = 0
$ws.Range("R29").Value = ''
$ws.Range("C29").Value = 46066

# New row 30  <=  old row 30  (A 17263-2022)
$ws.Range("A30:Z30").ClearContents()
$ws.Range("A30").Value = 'A 17263-2022'
$ws.Range("B30").Value = 44678
$ws.Range("D30").Value = 'STOCKHOLMS LÄN'
$ws.Range("E30").Value = 'ÖSTERÅKER'
$ws.Range("G30").Value = 2.4
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 0
$ws.Range("N30").Value = 0
$ws.Range("O30").Value = 0
$ws.Range("P30").Value = 0
$ws.Range("Q30").Value = 0
$ws.Range("R30").Value = ''
$ws.Range("C30").Value = 46066

# New row 31  <=  old row 31  (A 31075-2021)
$ws.Range("A31:Z31").ClearContents()
$ws.Range("A31").Value = 'A 31075-2021'
$ws.Range("B31").Value = 44368
$ws.Range("D31").Value = 'STOCKHOLMS LÄN'
$ws.Range("E31").Value = 'ÖSTERÅKER'
$ws.Range("F31").Value = 'Övriga Aktiebolag'
$ws.Range("G31").Value = 2.9
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = 0
$ws.Range("N31").Value = 0
$ws.Range("O31").Value = 0
$ws.Range("P31").Value = 0
$ws.Range("Q31").Value = 0
$ws.Range("R31").Value = ''
$ws.Range("C31").Value = 46066

# New row 32  <=  old row 32  (A 18172-2021)
$ws.Range("A32:Z32").ClearContents()
$ws.Range("A32").Value = 'A 18172-2021'
$ws.Range("B32").Value = 44302.62196759259
$ws.Range("D32").Value = 'STOCKHOLMS LÄN'
$ws.Range("E32").Value = 'ÖSTERÅKER'
$ws.Range("G32").Value = 3.2
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = 0
$ws.Range("N32").Value = 0
$ws.Range("O32").Value = 0
$ws.Range("P32").Value = 0
$ws.Range("Q32").Value = 0
$ws.Range("R32").Value = ''
$ws.Range("C32").Value = 46066

# New row 33  <=  old row 35  (A 67244-2021)
$ws.Range("A33:Z33").ClearContents()
$ws.Range("A33").Value = 'A 67244-2021'
$ws.Range("B33").Value = 44523
$ws.Range("D33").Value = 'STOCKHOLMS LÄN'
$ws.Range("E33").Value = 'ÖSTERÅKER'
$ws.Range("F33").Value = 'Kommuner'
$ws.Range("G33").Value = 5.3
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 0
$ws.Range("N33").Value = 0
$ws.Range("O33").Value = 0
$ws.Range("P33").Value = 0
$ws.Range("Q33").Value = 0
$ws.Range("R33").Value = ''
$ws.Range("C33").Value = 46066

# New row 34  <=  old row 40  (A 15131-2022)
$ws.Range("A34:Z34").ClearContents()
$ws.Range("A34").Value = 'A 15131-2022'
$ws.Range("B34").Value = 44658
$ws.Range("D34").Value = 'STOCKHOLMS LÄN'
$ws.Range("E34").Value = 'ÖSTERÅKER'
$ws.Range("G34").Value = 4.6
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = 0
$ws.Range("N34").Value = 0
$ws.Range("O34").Value = 0
$ws.Range("P34").Value = 0
$ws.Range("Q34").Value = 0
$ws.Range("R34").Value = ''
$ws.Range("C34").Value = 46066

# New row 35  <=  old row 42  (A 19934-2022)
$ws.Range("A35:Z35").ClearContents()
$ws.Range("A35").Value = 'A 19934-2022'
$ws.Range("B35").Value = 44697.46248842592
$ws.Range("D35").Value = 'STOCKHOLMS LÄN'
$ws.Range("E35").Value = 'ÖSTERÅKER'
$ws.Range("G35").Value = 0.3
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = 0
$ws.Range("N35").Value = 0
$ws.Range("O35").Value = 0
$ws.Range("P35").Value = 0
$ws.Range("Q35").Value = 0
$ws.Range("R35").Value = ''
$ws.Range("C35").Value = 46066

# New row 36  <=  old row 43  (A 19944-2022)
$ws.Range("A36:Z36").ClearContents()
$ws.Range("A36").Value = 'A 19944-2022'
$ws.Range("B36").Value = 44697.47424768518
$ws.Range("D36").Value = 'STOCKHOLMS LÄN'
$ws.Range("E36").Value = 'ÖSTERÅKER'
$ws.Range("G36").Value = 3.5
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = 0
$ws.Range("N36").Value = 0
$ws.Range("O36").Value = 0
$ws.Range("P36").Value = 0
$ws.Range("Q36").Value = 0
$ws.Range("R36").Value = ''
$ws.Range("C36").Value = 46066

# New row 37  <=  old row 37  (A 865-2022)
$ws.Range("A37:Z37").ClearContents()
$ws.Range("A37").Value = 'A 865-2022'
$ws.Range("B37").Value = 44569
$ws.Range("D37").Value = 'STOCKHOLMS LÄN'
$ws.Range("E37").Value = 'ÖSTERÅKER'
$ws.Range("F37").Value = 'Övriga Aktiebolag'
$ws.Range("G37").Value = 2.4
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = 0
$ws.Range("N37").Value = 0
$ws.Range("O37").Value = 0
$ws.Range("P37").Value = 0
$ws.Range("Q37").Value = 0
$ws.Range("R37").Value = ''
$ws.Range("C37").Value = 46066

# New row 38  <=  old row 36  (A 838-2022)
$ws.Range("A38:Z38").ClearContents()
$ws.Range("A38").Value = 'A 838-2022'
$ws.Range("B38").Value = 44569
$ws.Range("D38").Value = 'STOCKHOLMS LÄN'
$ws.Range("E38").Value = 'ÖSTERÅKER'
$ws.Range("F38").Value = 'Övriga Aktiebolag'
$ws.Range("G38").Value = 4.9
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 0
$ws.Range("N38").Value = 0
$ws.Range("O38").Value = 0
$ws.Range("P38").Value = 0
$ws.Range("Q38").Value = 0
$ws.Range("R38").Value = ''
$ws.Range("C38").Value = 46066

# New row 39  <=  old row 33  (A 13921-2021)
$ws.Range("A39:Z39").ClearContents()
$ws.Range("A39").Value = 'A 13921-2021'
$ws.Range("B39").Value = 44277
$ws.Range("D39").Value = 'STOCKHOLMS LÄN'
$ws.Range("E39").Value = 'ÖSTERÅKER'
$ws.Range("G39").Value = 0.9
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = 0
$ws.Range("N39").Value = 0
$ws.Range("O39").Value = 0
$ws.Range("P39").Value = 0
$ws.Range("Q39").Value = 0
$ws.Range("R39").Value = ''
$ws.Range("C39").Value = 46066

# New row 40  <=  old row 38  (A 28405-2022)
$ws.Range("A40:Z40").ClearContents()
$ws.Range("A40").Value = 'A 28405-2022'
$ws.Range("B40").Value = 44747
$ws.Range("D40").Value = 'STOCKHOLMS LÄN'
$ws.Range("E40").Value = 'ÖSTERÅKER'
$ws.Range("G40").Value = 1.9
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = 0
$ws.Range("N40").Value = 0
$ws.Range("O40").Value = 0
$ws.Range("P40").Value = 0
$ws.Range("Q40").Value = 0
$ws.Range("R40").Value = ''
$ws.Range("C40").Value = 46066

# New row 41  <=  old row 41  (A 17173-2022)
$ws.Range("A41:Z41").ClearContents()
$ws.Range("A41").Value = 'A 17173-2022'
$ws.Range("B41").Value = 44677
$ws.Range("D41").Value = 'STOCKHOLMS LÄN'
$ws.Range("E41").Value = 'ÖSTERÅKER'
$ws.Range("G41").Value = 2.8
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = 0
$ws.Range("N41").Value = 0
$ws.Range("O41").Value = 0
$ws.Range("P41").Value = 0
$ws.Range("Q41").Value = 0
$ws.Range("R41").Value = ''
$ws.Range("C41").Value = 46066

# New row 42  <=  old row 34  (A 10082-2022)
$ws.Range("A42:Z42").ClearContents()
$ws.Range("A42").Value = 'A 10082-2022'
$ws.Range("B42").Value = 44621
$ws.Range("D42").Value = 'STOCKHOLMS LÄN'
$ws.Range("E42").Value = 'ÖSTERÅKER'
$ws.Range("G42").Value = 6
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 0
$ws.Range("N42").Value = 0
$ws.Range("O42").Value = 0
$ws.Range("P42").Value = 0
$ws.Range("Q42").Value = 0
$ws.Range("R42").Value = ''
$ws.Range("C42").Value = 46066

# New row 43  <=  old row 39  (A 28964-2022)
$ws.Range("A43:Z43").ClearContents()
$ws.Range("A43").Value = 'A 28964-2022'
$ws.Range("B43").Value = 44749
$ws.Range("D43").Value = 'STOCKHOLMS LÄN'
$ws.Range("E43").Value = 'ÖSTERÅKER'
$ws.Range("G43").Value = 1.3
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = 0
$ws.Range("N43").Value = 0
$ws.Range("O43").Value = 0
$ws.Range("P43").Value = 0
$ws.Range("Q43").Value = 0
$ws.Range("R43").Value = ''
$ws.Range("C43").Value = 46066

# New row 44  <=  old row 44  (A 18091-2022)
$ws.Range("A44:Z44").ClearContents()
$ws.Range("A44").Value = 'A 18091-2022'
$ws.Range("B44").Value = 44684
$ws.Range("D44").Value = 'STOCKHOLMS LÄN'
$ws.Range("E44").Value = 'ÖSTERÅKER'
$ws.Range("G44").Value = 6.1
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = 0
$ws.Range("N44").Value = 0
$ws.Range("O44").Value = 0
$ws.Range("P44").Value = 0
$ws.Range("Q44").Value = 0
$ws.Range("R44").Value = ''
$ws.Range("C44").Value = 46066

# New row 45  <=  old row 46  (A 28399-2022)
$ws.Range("A45:Z45").ClearContents()
$ws.Range("A45").Value = 'A 28399-2022'
$ws.Range("B45").Value = 44747
$ws.Range("D45").Value = 'STOCKHOLMS LÄN'
$ws.Range("E45").Value = 'ÖSTERÅKER'
$ws.Range("G45").Value = 1.3
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = 0
$ws.Range("N45").Value = 0
$ws.Range("O45").Value = 0
$ws.Range("P45").Value = 0
$ws.Range("Q45").Value = 0
$ws.Range("R45").Value = ''
$ws.Range("C45").Value = 46066

# New row 46  <=  old row 95  (A 61370-2022)
$ws.Range("A46:Z46").ClearContents()
$ws.Range("A46").Value = 'A 61370-2022'
$ws.Range("B46").Value = 44915
$ws.Range("D46").Value = 'STOCKHOLMS LÄN'
$ws.Range("E46").Value = 'ÖSTERÅKER'
$ws.Range("G46").Value = 3.2
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = 0
$ws.Range("N46").Value = 0
$ws.Range("O46").Value = 0
$ws.Range("P46").Value = 0
$ws.Range("Q46").Value = 0
$ws.Range("R46").Value = ''
$ws.Range("C46").Value = 46066

# New row 47  <=  old row 91  (A 1603-2023)
$ws.Range("A47:Z47").ClearContents()
$ws.Range("A47").Value = 'A 1603-2023'
$ws.Range("B47").Value = 44937.68086805556
$ws.Range("D47").Value = 'STOCKHOLMS LÄN'
$ws.Range("E47").Value = 'ÖSTERÅKER'
$ws.Range("G47").Value = 1.7
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = 0
$ws.Range("N47").Value = 0
$ws.Range("O47").Value = 0
$ws.Range("P47").Value = 0
$ws.Range("Q47").Value = 0
$ws.Range("R47").Value = ''
$ws.Range("C47").Value = 46066

# New row 48  <=  old row 92  (A 1605-2023)
$ws.Range("A48:Z48").ClearContents()
$ws.Range("A48").Value = 'A 1605-2023'
$ws.Range("B48").Value = 44937.68299768519
$ws.Range("D48").Value = 'STOCKHOLMS LÄN'
$ws.Range("E48").Value = 'ÖSTERÅKER'
$ws.Range("G48").Value = 2.9
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = 0
$ws.Range("N48").Value = 0
$ws.Range("O48").Value = 0
$ws.Range("P48").Value = 0
$ws.Range("Q48").Value = 0
$ws.Range("R48").Value = ''
$ws.Range("C48").Value = 46066

# New row 49  <=  old row 104  (A 37037-2023)
$ws.Range("A49:Z49").ClearContents()
$ws.Range("A49").Value = 'A 37037-2023'
$ws.Range("B49").Value = 45155
$ws.Range("D49").Value = 'STOCKHOLMS LÄN'
$ws.Range("E49").Value = 'ÖSTERÅKER'
$ws.Range("G49").Value = 2
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = 0
$ws.Range("N49").Value = 0
$ws.Range("O49").Value = 0
$ws.Range("P49").Value = 0
$ws.Range("Q49").Value = 0
$ws.Range("R49").Value = ''
$ws.Range("C49").Value = 46066

# New row 50  <=  old row 97  (A 3828-2023)
$ws.Range("A50:Z50").ClearContents()
$ws.Range("A50").Value = 'A 3828-2023'
$ws.Range("B50").Value = 44951
$ws.Range("D50").Value = 'STOCKHOLMS LÄN'
$ws.Range("E50").Value = 'ÖSTERÅKER'
$ws.Range("F50").Value = 'Kommuner'
$ws.Range("G50").Value = 3.2
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = 0
$ws.Range("N50").Value = 0
$ws.Range("O50").Value = 0
$ws.Range("P50").Value = 0
$ws.Range("Q50").Value = 0
$ws.Range("R50").Value = ''
$ws.Range("C50").Value = 46066

# New row 51  <=  old row 51  (A 33158-2024)
$ws.Range("A51:Z51").ClearContents()
$ws.Range("A51").Value = 'A 33158-2024'
$ws.Range("B51").Value = 45518
$ws.Range("D51").Value = 'STOCKHOLMS LÄN'
$ws.Range("E51").Value = 'ÖSTERÅKER'
$ws.Range("G51").Value = 1.6
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = 0
$ws.Range("N51").Value = 0
$ws.Range("O51").Value = 0
$ws.Range("P51").Value = 0
$ws.Range("Q51").Value = 0
$ws.Range("R51").Value = ''
$ws.Range("C51").Value = 46066

# New row 52  <=  old row 114  (A 62831-2025)
$ws.Range("A52:Z52").ClearContents()
$ws.Range("A52").Value = 'A 62831-2025'
$ws.Range("B52").Value = 46008
$ws.Range("D52").Value = 'STOCKHOLMS LÄN'
$ws.Range("E52").Value = 'ÖSTERÅKER'
$ws.Range("G52").Value = 1.2
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = 0
$ws.Range("N52").Value = 0
$ws.Range("O52").Value = 0
$ws.Range("P52").Value = 0
$ws.Range("Q52").Value = 0
$ws.Range("R52").Value = ''
$ws.Range("C52").Value = 46066

# New row 53  <=  old row 115  (A 62839-2025)
$ws.Range("A53:Z53").ClearContents()
$ws.Range("A53").Value = 'A 62839-2025'
$ws.Range("B53").Value = 46008
$ws.Range("D53").Value = 'STOCKHOLMS LÄN'
$ws.Range("E53").Value = 'ÖSTERÅKER'
$ws.Range("G53").Value = 1.4
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = 0
$ws.Range("N53").Value = 0
$ws.Range("O53").Value = 0
$ws.Range("P53").Value = 0
$ws.Range("Q53").Value = 0
$ws.Range("R53").Value = ''
$ws.Range("C53").Value = 46066

# New row 54  <=  old row 49  (A 9865-2025)
$ws.Range("A54:Z54").ClearContents()
$ws.Range("A54").Value = 'A 9865-2025'
$ws.Range("B54").Value = 45716
$ws.Range("D54").Value = 'STOCKHOLMS LÄN'
$ws.Range("E54").Value = 'ÖSTERÅKER'
$ws.Range("G54").Value = 7
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = 0
$ws.Range("N54").Value = 0
$ws.Range("O54").Value = 0
$ws.Range("P54").Value = 0
$ws.Range("Q54").Value = 0
$ws.Range("R54").Value = ''
$ws.Range("C54").Value = 46066

# New row 55  <=  old row 62  (A 27137-2025)
$ws.Range("A55:Z55").ClearContents()
$ws.Range("A55").Value = 'A 27137-2025'
$ws.Range("B55").Value = 45811.95890046296
$ws.Range("D55").Value = 'STOCKHOLMS LÄN'
$ws.Range("E55").Value = 'ÖSTERÅKER'
$ws.Range("G55").Value = 3.8
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = 0
$ws.Range("N55").Value = 0
$ws.Range("O55").Value = 0
$ws.Range("P55").Value = 0
$ws.Range("Q55").Value = 0
$ws.Range("R55").Value = ''
$ws.Range("C55").Value = 46066

# New row 56  <=  old row 64  (A 27948-2025)
$ws.Range("A56:Z56").ClearContents()
$ws.Range("A56").Value = 'A 27948-2025'
$ws.Range("B56").Value = 45817.52460648148
$ws.Range("D56").Value = 'STOCKHOLMS LÄN'
$ws.Range("E56").Value = 'ÖSTERÅKER'
$ws.Range("G56").Value = 5.6
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = 0
$ws.Range("N56").Value = 0
$ws.Range("O56").Value = 0
$ws.Range("P56").Value = 0
$ws.Range("Q56").Value = 0
$ws.Range("R56").Value = ''
$ws.Range("C56").Value = 46066

# New row 57  <=  old row 65  (A 464-2022)
$ws.Range("A57:Z57").ClearContents()
$ws.Range("A57").Value = 'A 464-2022'
$ws.Range("B57").Value = 44565
$ws.Range("D57").Value = 'STOCKHOLMS LÄN'
$ws.Range("E57").Value = 'ÖSTERÅKER'
$ws.Range("F57").Value = 'Övriga Aktiebolag'
$ws.Range("G57").Value = 7.6
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = 0
$ws.Range("N57").Value = 0
$ws.Range("O57").Value = 0
$ws.Range("P57").Value = 0
$ws.Range("Q57").Value = 0
$ws.Range("R57").Value = ''
$ws.Range("C57").Value = 46066

# New row 58  <=  old row 63  (A 17353-2025)
$ws.Range("A58:Z58").ClearContents()
$ws.Range("A58").Value = 'A 17353-2025'
$ws.Range("B58").Value = 45756
$ws.Range("D58").Value = 'STOCKHOLMS LÄN'
$ws.Range("E58").Value = 'ÖSTERÅKER'
$ws.Range("F58").Value = 'Övriga Aktiebolag'
$ws.Range("G58").Value = 5.5
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = 0
$ws.Range("N58").Value = 0
$ws.Range("O58").Value = 0
$ws.Range("P58").Value = 0
$ws.Range("Q58").Value = 0
$ws.Range("R58").Value = ''
$ws.Range("C58").Value = 46066

# New row 59  <=  old row 53  (A 61360-2022)
$ws.Range("A59:Z59").ClearContents()
$ws.Range("A59").Value = 'A 61360-2022'
$ws.Range("B59").Value = 44915
$ws.Range("D59").Value = 'STOCKHOLMS LÄN'
$ws.Range("E59").Value = 'ÖSTERÅKER'
$ws.Range("G59").Value = 2.3
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = 0
$ws.Range("N59").Value = 0
$ws.Range("O59").Value = 0
$ws.Range("P59").Value = 0
$ws.Range("Q59").Value = 0
$ws.Range("R59").Value = ''
$ws.Range("C59").Value = 46066

# New row 60  <=  old row 71  (A 31456-2025)
$ws.Range("A60:Z60").ClearContents()
$ws.Range("A60").Value = 'A 31456-2025'
$ws.Range("B60").Value = 45833.5402662037
$ws.Range("D60").Value = 'STOCKHOLMS LÄN'
$ws.Range("E60").Value = 'ÖSTERÅKER'
$ws.Range("G60").Value = 4.3
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = 0
$ws.Range("N60").Value = 0
$ws.Range("O60").Value = 0
$ws.Range("P60").Value = 0
$ws.Range("Q60").Value = 0
$ws.Range("R60").Value = ''
$ws.Range("C60").Value = 46066

# New row 61  <=  old row 50  (A 61323-2022)
$ws.Range("A61:Z61").ClearContents()
$ws.Range("A61").Value = 'A 61323-2022'
$ws.Range("B61").Value = 44915
$ws.Range("D61").Value = 'STOCKHOLMS LÄN'
$ws.Range("E61").Value = 'ÖSTERÅKER'
$ws.Range("G61").Value = 2.7
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = 0
$ws.Range("N61").Value = 0
$ws.Range("O61").Value = 0
$ws.Range("P61").Value = 0
$ws.Range("Q61").Value = 0
$ws.Range("R61").Value = ''
$ws.Range("C61").Value = 46066

# New row 62  <=  old row 45  (A 22080-2024)
$ws.Range("A62:Z62").ClearContents()
$ws.Range("A62").Value = 'A 22080-2024'
$ws.Range("B62").Value = 45443.64962962963
$ws.Range("D62").Value = 'STOCKHOLMS LÄN'
$ws.Range("E62").Value = 'ÖSTERÅKER'
$ws.Range("F62").Value = 'Övriga Aktiebolag'
$ws.Range("G62").Value = 2.1
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = 0
$ws.Range("N62").Value = 0
$ws.Range("O62").Value = 0
$ws.Range("P62").Value = 0
$ws.Range("Q62").Value = 0
$ws.Range("R62").Value = ''
$ws.Range("C62").Value = 46066

# New row 63  <=  old row 89  (A 16728-2022)
$ws.Range("A63:Z63").ClearContents()
$ws.Range("A63").Value = 'A 16728-2022'
$ws.Range("B63").Value = 44673
$ws.Range("D63").Value = 'STOCKHOLMS LÄN'
$ws.Range("E63").Value = 'ÖSTERÅKER'
$ws.Range("G63").Value = 0.9
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = 0
$ws.Range("N63").Value = 0
$ws.Range("O63").Value = 0
$ws.Range("P63").Value = 0
$ws.Range("Q63").Value = 0
$ws.Range("R63").Value = ''
$ws.Range("C63").Value = 46066

# New row 64  <=  old row 78  (A 34038-2025)
$ws.Range("A64:Z64").ClearContents()
$ws.Range("A64").Value = 'A 34038-2025'
$ws.Range("B64").Value = 45842
$ws.Range("D64").Value = 'STOCKHOLMS LÄN'
$ws.Range("E64").Value = 'ÖSTERÅKER'
$ws.Range("F64").Value = 'Övriga Aktiebolag'
$ws.Range("G64").Value = 4.7
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = 0
$ws.Range("N64").Value = 0
$ws.Range("O64").Value = 0
$ws.Range("P64").Value = 0
$ws.Range("Q64").Value = 0
$ws.Range("R64").Value = ''
$ws.Range("C64").Value = 46066

# New row 65  <=  old row 79  (A 34079-2025)
$ws.Range("A65:Z65").ClearContents()
$ws.Range("A65").Value = 'A 34079-2025'
$ws.Range("B65").Value = 45842
$ws.Range("D65").Value = 'STOCKHOLMS LÄN'
$ws.Range("E65").Value = 'ÖSTERÅKER'
$ws.Range("F65").Value = 'Övriga Aktiebolag'
$ws.Range("G65").Value = 4.1
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = 0
$ws.Range("N65").Value = 0
$ws.Range("O65").Value = 0
$ws.Range("P65").Value = 0
$ws.Range("Q65").Value = 0
$ws.Range("R65").Value = ''
$ws.Range("C65").Value = 46066

# New row 66  <=  old row 76  (A 34034-2025)
$ws.Range("A66:Z66").ClearContents()
$ws.Range("A66").Value = 'A 34034-2025'
$ws.Range("B66").Value = 45842
$ws.Range("D66").Value = 'STOCKHOLMS LÄN'
$ws.Range("E66").Value = 'ÖSTERÅKER'
$ws.Range("F66").Value = 'Övriga Aktiebolag'
$ws.Range("G66").Value = 3.9
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = 0
$ws.Range("N66").Value = 0
$ws.Range("O66").Value = 0
$ws.Range("P66").Value = 0
$ws.Range("Q66").Value = 0
$ws.Range("R66").Value = ''
$ws.Range("C66").Value = 46066

# New row 67  <=  old row 100  (A 1136-2023)
$ws.Range("A67:Z67").ClearContents()
$ws.Range("A67").Value = 'A 1136-2023'
$ws.Range("B67").Value = 44935
$ws.Range("D67").Value = 'STOCKHOLMS LÄN'
$ws.Range("E67").Value = 'ÖSTERÅKER'
$ws.Range("G67").Value = 0.5
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = 0
$ws.Range("N67").Value = 0
$ws.Range("O67").Value = 0
$ws.Range("P67").Value = 0
$ws.Range("Q67").Value = 0
$ws.Range("R67").Value = ''
$ws.Range("C67").Value = 46066

# New row 68  <=  old row 61  (A 61710-2022)
$ws.Range("A68:Z68").ClearContents()
$ws.Range("A68").Value = 'A 61710-2022'
$ws.Range("B68").Value = 44917
$ws.Range("D68").Value = 'STOCKHOLMS LÄN'
$ws.Range("E68").Value = 'ÖSTERÅKER'
$ws.Range("G68").Value = 6.2
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = 0
$ws.Range("N68").Value = 0
$ws.Range("O68").Value = 0
$ws.Range("P68").Value = 0
$ws.Range("Q68").Value = 0
$ws.Range("R68").Value = ''
$ws.Range("C68").Value = 46066

# New row 69  <=  old row 57  (A 22528-2024)
$ws.Range("A69:Z69").ClearContents()
$ws.Range("A69").Value = 'A 22528-2024'
$ws.Range("B69").Value = 45447.44967592593
$ws.Range("D69").Value = 'STOCKHOLMS LÄN'
$ws.Range("E69").Value = 'ÖSTERÅKER'
$ws.Range("F69").Value = 'Övriga Aktiebolag'
$ws.Range("G69").Value = 1.4
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = 0
$ws.Range("N69").Value = 0
$ws.Range("O69").Value = 0
$ws.Range("P69").Value = 0
$ws.Range("Q69").Value = 0
$ws.Range("R69").Value = ''
$ws.Range("C69").Value = 46066

# New row 70  <=  old row 81  (A 35656-2025)
$ws.Range("A70:Z70").ClearContents()
$ws.Range("A70").Value = 'A 35656-2025'
$ws.Range("B70").Value = 45859
$ws.Range("D70").Value = 'STOCKHOLMS LÄN'
$ws.Range("E70").Value = 'ÖSTERÅKER'
$ws.Range("G70").Value = 1.5
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = 0
$ws.Range("N70").Value = 0
$ws.Range("O70").Value = 0
$ws.Range("P70").Value = 0
$ws.Range("Q70").Value = 0
$ws.Range("R70").Value = ''
$ws.Range("C70").Value = 46066

# New row 71  <=  old row 82  (A 35819-2025)
$ws.Range("A71:Z71").ClearContents()
$ws.Range("A71").Value = 'A 35819-2025'
$ws.Range("B71").Value = 45861
$ws.Range("D71").Value = 'STOCKHOLMS LÄN'
$ws.Range("E71").Value = 'ÖSTERÅKER'
$ws.Range("F71").Value = 'Övriga Aktiebolag'
$ws.Range("G71").Value = 8
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = 0
$ws.Range("N71").Value = 0
$ws.Range("O71").Value = 0
$ws.Range("P71").Value = 0
$ws.Range("Q71").Value = 0
$ws.Range("R71").Value = ''
$ws.Range("C71").Value = 46066

# New row 72  <=  old row 87  (A 36611-2025)
$ws.Range("A72:Z72").ClearContents()
$ws.Range("A72").Value = 'A 36611-2025'
$ws.Range("B72").Value = 45870
$ws.Range("D72").Value = 'STOCKHOLMS LÄN'
$ws.Range("E72").Value = 'ÖSTERÅKER'
$ws.Range("F72").Value = 'Övriga Aktiebolag'
$ws.Range("G72").Value = 2.6
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = 0
$ws.Range("N72").Value = 0
$ws.Range("O72").Value = 0
$ws.Range("P72").Value = 0
$ws.Range("Q72").Value = 0
$ws.Range("R72").Value = ''
$ws.Range("C72").Value = 46066

# New row 73  <=  old row 47  (A 28143-2024)
$ws.Range("A73:Z73").ClearContents()
$ws.Range("A73").Value = 'A 28143-2024'
$ws.Range("B73").Value = 45476
$ws.Range("D73").Value = 'STOCKHOLMS LÄN'
$ws.Range("E73").Value = 'ÖSTERÅKER'
$ws.Range("F73").Value = 'Övriga Aktiebolag'
$ws.Range("G73").Value = 2.2
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = 0
$ws.Range("N73").Value = 0
$ws.Range("O73").Value = 0
$ws.Range("P73").Value = 0
$ws.Range("Q73").Value = 0
$ws.Range("R73").Value = ''
$ws.Range("C73").Value = 46066

# New row 74  <=  old row 84  (A 36614-2025)
$ws.Range("A74:Z74").ClearContents()
$ws.Range("A74").Value = 'A 36614-2025'
$ws.Range("B74").Value = 45870
$ws.Range("D74").Value = 'STOCKHOLMS LÄN'
$ws.Range("E74").Value = 'ÖSTERÅKER'
$ws.Range("G74").Value = 1.7
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 0
$ws.Range("N74").Value = 0
$ws.Range("O74").Value = 0
$ws.Range("P74").Value = 0
$ws.Range("Q74").Value = 0
$ws.Range("R74").Value = ''
$ws.Range("C74").Value = 46066

# New row 75  <=  old row 86  (A 36612-2025)
$ws.Range("A75:Z75").ClearContents()
$ws.Range("A75").Value = 'A 36612-2025'
$ws.Range("B75").Value = 45870
$ws.Range("D75").Value = 'STOCKHOLMS LÄN'
$ws.Range("E75").Value = 'ÖSTERÅKER'
$ws.Range("F75").Value = 'Övriga Aktiebolag'
$ws.Range("G75").Value = 2.1
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = 0
$ws.Range("N75").Value = 0
$ws.Range("O75").Value = 0
$ws.Range("P75").Value = 0
$ws.Range("Q75").Value = 0
$ws.Range("R75").Value = ''
$ws.Range("C75").Value = 46066

# New row 76  <=  old row 106  (A 34713-2025)
$ws.Range("A76:Z76").ClearContents()
$ws.Range("A76").Value = 'A 34713-2025'
$ws.Range("B76").Value = 45848
$ws.Range("D76").Value = 'STOCKHOLMS LÄN'
$ws.Range("E76").Value = 'ÖSTERÅKER'
$ws.Range("G76").Value = 0.9
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = 0
$ws.Range("N76").Value = 0
$ws.Range("O76").Value = 0
$ws.Range("P76").Value = 0
$ws.Range("Q76").Value = 0
$ws.Range("R76").Value = ''
$ws.Range("C76").Value = 46066

# New row 77  <=  old row 90  (A 14776-2024)
$ws.Range("A77:Z77").ClearContents()
$ws.Range("A77").Value = 'A 14776-2024'
$ws.Range("B77").Value = 45397
$ws.Range("D77").Value = 'STOCKHOLMS LÄN'
$ws.Range("E77").Value = 'ÖSTERÅKER'
$ws.Range("G77").Value = 13.6
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 0
$ws.Range("N77").Value = 0
$ws.Range("O77").Value = 0
$ws.Range("P77").Value = 0
$ws.Range("Q77").Value = 0
$ws.Range("R77").Value = ''
$ws.Range("C77").Value = 46066

# New row 78  <=  old row 101  (A 39394-2023)
$ws.Range("A78:Z78").ClearContents()
$ws.Range("A78").Value = 'A 39394-2023'
$ws.Range("B78").Value = 45166
$ws.Range("D78").Value = 'STOCKHOLMS LÄN'
$ws.Range("E78").Value = 'ÖSTERÅKER'
$ws.Range("G78").Value = 0.4
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = 0
$ws.Range("N78").Value = 0
$ws.Range("O78").Value = 0
$ws.Range("P78").Value = 0
$ws.Range("Q78").Value = 0
$ws.Range("R78").Value = ''
$ws.Range("C78").Value = 46066

# New row 79  <=  old row 58  (A 16803-2022)
$ws.Range("A79:Z79").ClearContents()
$ws.Range("A79").Value = 'A 16803-2022'
$ws.Range("B79").Value = 44673
$ws.Range("D79").Value = 'STOCKHOLMS LÄN'
$ws.Range("E79").Value = 'ÖSTERÅKER'
$ws.Range("G79").Value = 1.9
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = 0
$ws.Range("N79").Value = 0
$ws.Range("O79").Value = 0
$ws.Range("P79").Value = 0
$ws.Range("Q79").Value = 0
$ws.Range("R79").Value = ''
$ws.Range("C79").Value = 46066

# New row 80  <=  old row 93  (A 470-2022)
$ws.Range("A80:Z80").ClearContents()
$ws.Range("A80").Value = 'A 470-2022'
$ws.Range("B80").Value = 44565
$ws.Range("D80").Value = 'STOCKHOLMS LÄN'
$ws.Range("E80").Value = 'ÖSTERÅKER'
$ws.Range("F80").Value = 'Övriga Aktiebolag'
$ws.Range("G80").Value = 4.6
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = 0
$ws.Range("N80").Value = 0
$ws.Range("O80").Value = 0
$ws.Range("P80").Value = 0
$ws.Range("Q80").Value = 0
$ws.Range("R80").Value = ''
$ws.Range("C80").Value = 46066

# New row 81  <=  old row 70  (A 23144-2024)
$ws.Range("A81:Z81").ClearContents()
$ws.Range("A81").Value = 'A 23144-2024'
$ws.Range("B81").Value = 45450.63202546296
$ws.Range("D81").Value = 'STOCKHOLMS LÄN'
$ws.Range("E81").Value = 'ÖSTERÅKER'
$ws.Range("F81").Value = 'Övriga Aktiebolag'
$ws.Range("G81").Value = 4.7
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = 0
$ws.Range("N81").Value = 0
$ws.Range("O81").Value = 0
$ws.Range("P81").Value = 0
$ws.Range("Q81").Value = 0
$ws.Range("R81").Value = ''
$ws.Range("C81").Value = 46066

# New row 82  <=  old row 55  (A 13118-2021)
$ws.Range("A82:Z82").ClearContents()
$ws.Range("A82").Value = 'A 13118-2021'
$ws.Range("B82").Value = 44272
$ws.Range("D82").Value = 'STOCKHOLMS LÄN'
$ws.Range("E82").Value = 'ÖSTERÅKER'
$ws.Range("F82").Value = 'Kyrkan'
$ws.Range("G82").Value = 4
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = 0
$ws.Range("N82").Value = 0
$ws.Range("O82").Value = 0
$ws.Range("P82").Value = 0
$ws.Range("Q82").Value = 0
$ws.Range("R82").Value = ''
$ws.Range("C82").Value = 46066

# New row 83  <=  old row 110  (A 52220-2025)
$ws.Range("A83:Z83").ClearContents()
$ws.Range("A83").Value = 'A 52220-2025'
$ws.Range("B83").Value = 45953
$ws.Range("D83").Value = 'STOCKHOLMS LÄN'
$ws.Range("E83").Value = 'ÖSTERÅKER'
$ws.Range("G83").Value = 0.8
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = 0
$ws.Range("N83").Value = 0
$ws.Range("O83").Value = 0
$ws.Range("P83").Value = 0
$ws.Range("Q83").Value = 0
$ws.Range("R83").Value = ''
$ws.Range("C83").Value = 46066

# New row 84  <=  old row 105  (A 19404-2024)
$ws.Range("A84:Z84").ClearContents()
$ws.Range("A84").Value = 'A 19404-2024'
$ws.Range("B84").Value = 45428
$ws.Range("D84").Value = 'STOCKHOLMS LÄN'
$ws.Range("E84").Value = 'ÖSTERÅKER'
$ws.Range("F84").Value = 'Övriga Aktiebolag'
$ws.Range("G84").Value = 2.2
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 0
$ws.Range("N84").Value = 0
$ws.Range("O84").Value = 0
$ws.Range("P84").Value = 0
$ws.Range("Q84").Value = 0
$ws.Range("R84").Value = ''
$ws.Range("C84").Value = 46066

# New row 85  <=  old row 60  (A 6655-2023)
$ws.Range("A85:Z85").ClearContents()
$ws.Range("A85").Value = 'A 6655-2023'
$ws.Range("B85").Value = 44966.60668981481
$ws.Range("D85").Value = 'STOCKHOLMS LÄN'
$ws.Range("E85").Value = 'ÖSTERÅKER'
$ws.Range("G85").Value = 0.8
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 0
$ws.Range("N85").Value = 0
$ws.Range("O85").Value = 0
$ws.Range("P85").Value = 0
$ws.Range("Q85").Value = 0
$ws.Range("R85").Value = ''
$ws.Range("C85").Value = 46066

# New row 86  <=  old row 74  (A 9842-2025)
$ws.Range("A86:Z86").ClearContents()
$ws.Range("A86").Value = 'A 9842-2025'
$ws.Range("B86").Value = 45716
$ws.Range("D86").Value = 'STOCKHOLMS LÄN'
$ws.Range("E86").Value = 'ÖSTERÅKER'
$ws.Range("G86").Value = 4.8
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 0
$ws.Range("N86").Value = 0
$ws.Range("O86").Value = 0
$ws.Range("P86").Value = 0
$ws.Range("Q86").Value = 0
$ws.Range("R86").Value = ''
$ws.Range("C86").Value = 46066

# New row 87  <=  old row 75  (A 56909-2024)
$ws.Range("A87:Z87").ClearContents()
$ws.Range("A87").Value = 'A 56909-2024'
$ws.Range("B87").Value = 45628.59060185185
$ws.Range("D87").Value = 'STOCKHOLMS LÄN'
$ws.Range("E87").Value = 'ÖSTERÅKER'
$ws.Range("G87").Value = 1
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = 0
$ws.Range("N87").Value = 0
$ws.Range("O87").Value = 0
$ws.Range("P87").Value = 0
$ws.Range("Q87").Value = 0
$ws.Range("R87").Value = ''
$ws.Range("C87").Value = 46066

# New row 88  <=  old row 102  (A 30949-2024)
$ws.Range("A88:Z88").ClearContents()
$ws.Range("A88").Value = 'A 30949-2024'
$ws.Range("B88").Value = 45499
$ws.Range("D88").Value = 'STOCKHOLMS LÄN'
$ws.Range("E88").Value = 'ÖSTERÅKER'
$ws.Range("F88").Value = 'Övriga Aktiebolag'
$ws.Range("G88").Value = 3.2
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = 0
$ws.Range("N88").Value = 0
$ws.Range("O88").Value = 0
$ws.Range("P88").Value = 0
$ws.Range("Q88").Value = 0
$ws.Range("R88").Value = ''
$ws.Range("C88").Value = 46066

# New row 89  <=  old row 103  (A 845-2024)
$ws.Range("A89:Z89").ClearContents()
$ws.Range("A89").Value = 'A 845-2024'
$ws.Range("B89").Value = 45300.91423611111
$ws.Range("D89").Value = 'STOCKHOLMS LÄN'
$ws.Range("E89").Value = 'ÖSTERÅKER'
$ws.Range("G89").Value = 10.1
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 0
$ws.Range("N89").Value = 0
$ws.Range("O89").Value = 0
$ws.Range("P89").Value = 0
$ws.Range("Q89").Value = 0
$ws.Range("R89").Value = ''
$ws.Range("C89").Value = 46066

# New row 90  <=  old row 66  (A 31134-2023)
$ws.Range("A90:Z90").ClearContents()
$ws.Range("A90").Value = 'A 31134-2023'
$ws.Range("B90").Value = 45113
$ws.Range("D90").Value = 'STOCKHOLMS LÄN'
$ws.Range("E90").Value = 'ÖSTERÅKER'
$ws.Range("G90").Value = 2.1
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = 0
$ws.Range("N90").Value = 0
$ws.Range("O90").Value = 0
$ws.Range("P90").Value = 0
$ws.Range("Q90").Value = 0
$ws.Range("R90").Value = ''
$ws.Range("C90").Value = 46066

# New row 91  <=  old row 48  (A 31575-2022)
$ws.Range("A91:Z91").ClearContents()
$ws.Range("A91").Value = 'A 31575-2022'
$ws.Range("B91").Value = 44775.45063657407
$ws.Range("D91").Value = 'STOCKHOLMS LÄN'
$ws.Range("E91").Value = 'ÖSTERÅKER'
$ws.Range("G91").Value = 5
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = 0
$ws.Range("N91").Value = 0
$ws.Range("O91").Value = 0
$ws.Range("P91").Value = 0
$ws.Range("Q91").Value = 0
$ws.Range("R91").Value = ''
$ws.Range("C91").Value = 46066

# New row 92  <=  old row 111  (A 60482-2025)
$ws.Range("A92:Z92").ClearContents()
$ws.Range("A92").Value = 'A 60482-2025'
$ws.Range("B92").Value = 45995.59832175926
$ws.Range("D92").Value = 'STOCKHOLMS LÄN'
$ws.Range("E92").Value = 'ÖSTERÅKER'
$ws.Range("G92").Value = 1.4
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 0
$ws.Range("N92").Value = 0
$ws.Range("O92").Value = 0
$ws.Range("P92").Value = 0
$ws.Range("Q92").Value = 0
$ws.Range("R92").Value = ''
$ws.Range("C92").Value = 46066

# New row 93  <=  old row 88  (A 10185-2025)
$ws.Range("A93:Z93").ClearContents()
$ws.Range("A93").Value = 'A 10185-2025'
$ws.Range("B93").Value = 45719
$ws.Range("D93").Value = 'STOCKHOLMS LÄN'
$ws.Range("E93").Value = 'ÖSTERÅKER'
$ws.Range("G93").Value = 4.1
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 0
$ws.Range("N93").Value = 0
$ws.Range("O93").Value = 0
$ws.Range("P93").Value = 0
$ws.Range("Q93").Value = 0
$ws.Range("R93").Value = ''
$ws.Range("C93").Value = 46066

# New row 94  <=  old row 112  (A 61781-2025)
$ws.Range("A94:Z94").ClearContents()
$ws.Range("A94").Value = 'A 61781-2025'
$ws.Range("B94").Value = 46002
$ws.Range("D94").Value = 'STOCKHOLMS LÄN'
$ws.Range("E94").Value = 'ÖSTERÅKER'
$ws.Range("G94").Value = 4.9
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = 0
$ws.Range("N94").Value = 0
$ws.Range("O94").Value = 0
$ws.Range("P94").Value = 0
$ws.Range("Q94").Value = 0
$ws.Range("R94").Value = ''
$ws.Range("C94").Value = 46066

# New row 95  <=  old row 113  (A 5038-2026)
$ws.Range("A95:Z95").ClearContents()
$ws.Range("A95").Value = 'A 5038-2026'
$ws.Range("B95").Value = 46049.46427083333
$ws.Range("D95").Value = 'STOCKHOLMS LÄN'
$ws.Range("E95").Value = 'ÖSTERÅKER'
$ws.Range("G95").Value = 2.1
$ws.Range("H95").Value = 0
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("M95").Value = 0
$ws.Range("N95").Value = 0
$ws.Range("O95").Value = 0
$ws.Range("P95").Value = 0
$ws.Range("Q95").Value = 0
$ws.Range("R95").Value = ''
$ws.Range("C95").Value = 46066

# New row 96  <=  old row 94  (A 67265-2021)
$ws.Range("A96:Z96").ClearContents()
$ws.Range("A96").Value = 'A 67265-2021'
$ws.Range("B96").Value = 44523.59277777778
$ws.Range("D96").Value = 'STOCKHOLMS LÄN'
$ws.Range("E96").Value = 'ÖSTERÅKER'
$ws.Range("F96").Value = 'Kommuner'
$ws.Range("G96").Value = 2.6
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = 0
$ws.Range("N96").Value = 0
$ws.Range("O96").Value = 0
$ws.Range("P96").Value = 0
$ws.Range("Q96").Value = 0
$ws.Range("R96").Value = ''
$ws.Range("C96").Value = 46066

# New row 97  <=  old row 99  (A 18711-2023)
$ws.Range("A97:Z97").ClearContents()
$ws.Range("A97").Value = 'A 18711-2023'
$ws.Range("B97").Value = 45043.6433912037
$ws.Range("D97").Value = 'STOCKHOLMS LÄN'
$ws.Range("E97").Value = 'ÖSTERÅKER'
$ws.Range("G97").Value = 3
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 0
$ws.Range("N97").Value = 0
$ws.Range("O97").Value = 0
$ws.Range("P97").Value = 0
$ws.Range("Q97").Value = 0
$ws.Range("R97").Value = ''
$ws.Range("C97").Value = 46066

# New row 98  <=  old row 72  (A 24499-2023)
$ws.Range("A98:Z98").ClearContents()
$ws.Range("A98").Value = 'A 24499-2023'
$ws.Range("B98").Value = 45082
$ws.Range("D98").Value = 'STOCKHOLMS LÄN'
$ws.Range("E98").Value = 'ÖSTERÅKER'
$ws.Range("G98").Value = 0.1
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 0
$ws.Range("N98").Value = 0
$ws.Range("O98").Value = 0
$ws.Range("P98").Value = 0
$ws.Range("Q98").Value = 0
$ws.Range("R98").Value = ''
$ws.Range("C98").Value = 46066

# New row 99  <=  old row 54  (A 28146-2024)
$ws.Range("A99:Z99").ClearContents()
$ws.Range("A99").Value = 'A 28146-2024'
$ws.Range("B99").Value = 45476
$ws.Range("D99").Value = 'STOCKHOLMS LÄN'
$ws.Range("E99").Value = 'ÖSTERÅKER'
$ws.Range("F99").Value = 'Övriga Aktiebolag'
$ws.Range("G99").Value = 2.7
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 0
$ws.Range("N99").Value = 0
$ws.Range("O99").Value = 0
$ws.Range("P99").Value = 0
$ws.Range("Q99").Value = 0
$ws.Range("R99").Value = ''
$ws.Range("C99").Value = 46066

# New row 100  <=  old row 68  (A 57048-2024)
$ws.Range("A100:Z100").ClearContents()
$ws.Range("A100").Value = 'A 57048-2024'
$ws.Range("B100").Value = 45629.33938657407
$ws.Range("D100").Value = 'STOCKHOLMS LÄN'
$ws.Range("E100").Value = 'ÖSTERÅKER'
$ws.Range("G100").Value = 3
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = 0
$ws.Range("N100").Value = 0
$ws.Range("O100").Value = 0
$ws.Range("P100").Value = 0
$ws.Range("Q100").Value = 0
$ws.Range("R100").Value = ''
$ws.Range("C100").Value = 46066

# New row 101  <=  old row 98  (A 5286-2025)
$ws.Range("A101:Z101").ClearContents()
$ws.Range("A101").Value = 'A 5286-2025'
$ws.Range("B101").Value = 45692.44144675926
$ws.Range("D101").Value = 'STOCKHOLMS LÄN'
$ws.Range("E101").Value = 'ÖSTERÅKER'
$ws.Range("G101").Value = 11.2
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = 0
$ws.Range("N101").Value = 0
$ws.Range("O101").Value = 0
$ws.Range("P101").Value = 0
$ws.Range("Q101").Value = 0
$ws.Range("R101").Value = ''
$ws.Range("C101").Value = 46066

# New row 102  <=  old row 96  (A 3817-2023)
$ws.Range("A102:Z102").ClearContents()
$ws.Range("A102").Value = 'A 3817-2023'
$ws.Range("B102").Value = 44951
$ws.Range("D102").Value = 'STOCKHOLMS LÄN'
$ws.Range("E102").Value = 'ÖSTERÅKER'
$ws.Range("F102").Value = 'Kommuner'
$ws.Range("G102").Value = 4.7
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 0
$ws.Range("N102").Value = 0
$ws.Range("O102").Value = 0
$ws.Range("P102").Value = 0
$ws.Range("Q102").Value = 0
$ws.Range("R102").Value = ''
$ws.Range("C102").Value = 46066

# New row 103  <=  old row 83  (A 544-2023)
$ws.Range("A103:Z103").ClearContents()
$ws.Range("A103").Value = 'A 544-2023'
$ws.Range("B103").Value = 44924
$ws.Range("D103").Value = 'STOCKHOLMS LÄN'
$ws.Range("E103").Value = 'ÖSTERÅKER'
$ws.Range("F103").Value = 'Övriga Aktiebolag'
$ws.Range("G103").Value = 0.9
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = 0
$ws.Range("N103").Value = 0
$ws.Range("O103").Value = 0
$ws.Range("P103").Value = 0
$ws.Range("Q103").Value = 0
$ws.Range("R103").Value = ''
$ws.Range("C103").Value = 46066

# New row 104  <=  old row 73  (A 3326-2022)
$ws.Range("A104:Z104").ClearContents()
$ws.Range("A104").Value = 'A 3326-2022'
$ws.Range("B104").Value = 44583
$ws.Range("D104").Value = 'STOCKHOLMS LÄN'
$ws.Range("E104").Value = 'ÖSTERÅKER'
$ws.Range("F104").Value = 'Övriga Aktiebolag'
$ws.Range("G104").Value = 2.7
$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("M104").Value = 0
$ws.Range("N104").Value = 0
$ws.Range("O104").Value = 0
$ws.Range("P104").Value = 0
$ws.Range("Q104").Value = 0
$ws.Range("R104").Value = ''
$ws.Range("C104").Value = 46066

# New row 105  <=  old row 52  (A 13366-2024)
$ws.Range("A105:Z105").ClearContents()
$ws.Range("A105").Value = 'A 13366-2024'
$ws.Range("B105").Value = 45387
$ws.Range("D105").Value = 'STOCKHOLMS LÄN'
$ws.Range("E105").Value = 'ÖSTERÅKER'
$ws.Range("G105").Value = 13.7
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 0
$ws.Range("N105").Value = 0
$ws.Range("O105").Value = 0
$ws.Range("P105").Value = 0
$ws.Range("Q105").Value = 0
$ws.Range("R105").Value = ''
$ws.Range("C105").Value = 46066

# New row 106  <=  old row 69  (A 57055-2024)
$ws.Range("A106:Z106").ClearContents()
$ws.Range("A106").Value = 'A 57055-2024'
$ws.Range("B106").Value = 45629.34627314815
$ws.Range("D106").Value = 'STOCKHOLMS LÄN'
$ws.Range("E106").Value = 'ÖSTERÅKER'
$ws.Range("G106").Value = 0.5
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = 0
$ws.Range("N106").Value = 0
$ws.Range("O106").Value = 0
$ws.Range("P106").Value = 0
$ws.Range("Q106").Value = 0
$ws.Range("R106").Value = ''
$ws.Range("C106").Value = 46066

# New row 107  <=  old row 85  (A 59379-2022)
$ws.Range("A107:Z107").ClearContents()
$ws.Range("A107").Value = 'A 59379-2022'
$ws.Range("B107").Value = 44905
$ws.Range("D107").Value = 'STOCKHOLMS LÄN'
$ws.Range("E107").Value = 'ÖSTERÅKER'
$ws.Range("F107").Value = 'Kommuner'
$ws.Range("G107").Value = 4.7
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 0
$ws.Range("N107").Value = 0
$ws.Range("O107").Value = 0
$ws.Range("P107").Value = 0
$ws.Range("Q107").Value = 0
$ws.Range("R107").Value = ''
$ws.Range("C107").Value = 46066

# New row 108  <=  old row 56  (A 46688-2022)
$ws.Range("A108:Z108").ClearContents()
$ws.Range("A108").Value = 'A 46688-2022'
$ws.Range("B108").Value = 44851
$ws.Range("D108").Value = 'STOCKHOLMS LÄN'
$ws.Range("E108").Value = 'ÖSTERÅKER'
$ws.Range("G108").Value = 2
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = 0
$ws.Range("N108").Value = 0
$ws.Range("O108").Value = 0
$ws.Range("P108").Value = 0
$ws.Range("Q108").Value = 0
$ws.Range("R108").Value = ''
$ws.Range("C108").Value = 46066

# New row 109  <=  old row 80  (A 16738-2022)
$ws.Range("A109:Z109").ClearContents()
$ws.Range("A109").Value = 'A 16738-2022'
$ws.Range("B109").Value = 44673
$ws.Range("D109").Value = 'STOCKHOLMS LÄN'
$ws.Range("E109").Value = 'ÖSTERÅKER'
$ws.Range("G109").Value = 1.1
$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = 0
$ws.Range("N109").Value = 0
$ws.Range("O109").Value = 0
$ws.Range("P109").Value = 0
$ws.Range("Q109").Value = 0
$ws.Range("R109").Value = ''
$ws.Range("C109").Value = 46066

# New row 110  <=  old row 77  (A 17196-2022)
$ws.Range("A110:Z110").ClearContents()
$ws.Range("A110").Value = 'A 17196-2022'
$ws.Range("B110").Value = 44677
$ws.Range("D110").Value = 'STOCKHOLMS LÄN'
$ws.Range("E110").Value = 'ÖSTERÅKER'
$ws.Range("G110").Value = 5.9
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 0
$ws.Range("N110").Value = 0
$ws.Range("O110").Value = 0
$ws.Range("P110").Value = 0
$ws.Range("Q110").Value = 0
$ws.Range("R110").Value = ''
$ws.Range("C110").Value = 46066

# New row 111  <=  old row 67  (A 3917-2023)
$ws.Range("A111:Z111").ClearContents()
$ws.Range("A111").Value = 'A 3917-2023'
$ws.Range("B111").Value = 44952
$ws.Range("D111").Value = 'STOCKHOLMS LÄN'
$ws.Range("E111").Value = 'ÖSTERÅKER'
$ws.Range("F111").Value = 'Kommuner'
$ws.Range("G111").Value = 2.8
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = 0
$ws.Range("N111").Value = 0
$ws.Range("O111").Value = 0
$ws.Range("P111").Value = 0
$ws.Range("Q111").Value = 0
$ws.Range("R111").Value = ''
$ws.Range("C111").Value = 46066

# New row 112  <=  old row 59  (A 62108-2022)
$ws.Range("A112:Z112").ClearContents()
$ws.Range("A112").Value = 'A 62108-2022'
$ws.Range("B112").Value = 44920
$ws.Range("D112").Value = 'STOCKHOLMS LÄN'
$ws.Range("E112").Value = 'ÖSTERÅKER'
$ws.Range("G112").Value = 1.5
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("M112").Value = 0
$ws.Range("N112").Value = 0
$ws.Range("O112").Value = 0
$ws.Range("P112").Value = 0
$ws.Range("Q112").Value = 0
$ws.Range("R112").Value = ''
$ws.Range("C112").Value = 46066

# New row 113  <=  old row 107  (A 37244-2024)
$ws.Range("A113:Z113").ClearContents()
$ws.Range("A113").Value = 'A 37244-2024'
$ws.Range("B113").Value = 45540
$ws.Range("D113").Value = 'STOCKHOLMS LÄN'
$ws.Range("E113").Value = 'ÖSTERÅKER'
$ws.Range("G113").Value = 2.6
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 0
$ws.Range("N113").Value = 0
$ws.Range("O113").Value = 0
$ws.Range("P113").Value = 0
$ws.Range("Q113").Value = 0
$ws.Range("R113").Value = ''
$ws.Range("C113").Value = 46066

# New row 114  <=  old row 109  (A 21796-2025)
$ws.Range("A114:Z114").ClearContents()
$ws.Range("A114").Value = 'A 21796-2025'
$ws.Range("B114").Value = 45783.68135416666
$ws.Range("D114").Value = 'STOCKHOLMS LÄN'
$ws.Range("E114").Value = 'ÖSTERÅKER'
$ws.Range("G114").Value = 2.5
$ws.Range("H114").Value = 0
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("M114").Value = 0
$ws.Range("N114").Value = 0
$ws.Range("O114").Value = 0
$ws.Range("P114").Value = 0
$ws.Range("Q114").Value = 0
$ws.Range("R114").Value = ''
$ws.Range("C114").Value = 46066

# New row 115  <=  old row 108  (A 21786-2025)
$ws.Range("A115:Z115").ClearContents()
$ws.Range("A115").Value = 'A 21786-2025'
$ws.Range("B115").Value = 45783.67215277778
$ws.Range("D115").Value = 'STOCKHOLMS LÄN'
$ws.Range("E115").Value = 'ÖSTERÅKER'
$ws.Range("G115").Value = 3.8
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = 0
$ws.Range("N115").Value = 0
$ws.Range("O115").Value = 0
$ws.Range("P115").Value = 0
$ws.Range("Q115").Value = 0
$ws.Range("R115").Value = ''
$ws.Range("C115").Value = 46066

